# edit.ps1 -- applies the "added final stuff blahj" commit to the deck:
#   1. Bumps the footer date placeholder text from 4/25/2021 to 4/28/2021
#      on the slide master and every slide layout (12 shapes total).
#   2. Updates the "4 oscillators" bullet on slide 2.
#   3. Updates the "Four Oscillator" bullet on slide 3 (5 oscillators, with
#      an extra ", Custom" inserted just before the closing paren).

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Date placeholder: slide master + all slide layouts.
# -----------------------------------------------------------------
$oldDate = "4/25/2021"
$newDate = "4/28/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# -----------------------------------------------------------------
# 2) Slide 2 - "Synthesizer has 4 oscillators ..." bullet.
# -----------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shp2 = $slide2.Shapes.Item("Content Placeholder 9")
$tr2 = $shp2.TextFrame.TextRange
$full2 = $tr2.Text
$oldLine2 = "Synthesizer has 4 oscillators (Sin, Sawtooth, Triangle, and Square)."
$newLine2 = "Synthesizer has 5 oscillators (Sin, Sawtooth, Triangle, Square, and Custom) with Octave changes."
$pos2 = $full2.IndexOf($oldLine2)
if ($pos2 -ge 0) {
    $range2 = $tr2.Characters($pos2 + 1, $oldLine2.Length)
    $range2.Text = $newLine2
}

# -----------------------------------------------------------------
# 3) Slide 3 - "Four Oscillator ( Sin, Sawtooth, Triangle, Square)" bullet.
#    Becomes "Five Oscillator ( Sin, Sawtooth, Triangle, Square, Custom)",
#    renaming Four->Five in place and then inserting ", Custom" right
#    before the trailing ")".
# -----------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shp3 = $slide3.Shapes.Item("Content Placeholder 2")
$tr3 = $shp3.TextFrame.TextRange
$full3 = $tr3.Text
$oldLine3 = "Four Oscillator ( Sin, Sawtooth, Triangle, Square)"
$newLine3 = "Five Oscillator ( Sin, Sawtooth, Triangle, Square)"
$pos3 = $full3.IndexOf($oldLine3)
if ($pos3 -ge 0) {
    # Replace the whole run's text in one shot (Four -> Five) so the run
    # stays a single run, keeping its original rPr/dirty state.
    $lineRange = $tr3.Characters($pos3 + 1, $oldLine3.Length)
    $lineRange.Text = $newLine3

    # Re-locate the (now renamed) line and insert ", Custom" before the
    # closing parenthesis, splitting the run the same way PowerPoint does
    # when new text is appended mid-paragraph.
    $full3b = $tr3.Text
    $pos3b = $full3b.IndexOf($newLine3)
    $closeParenPos = $pos3b + $newLine3.Length
    $parenRange = $tr3.Characters($closeParenPos, 1)
    [void]$parenRange.InsertBefore(", Custom")
}
